# Updated the data masking logic
# Translate the PII-related column headers on the "Student_Information"
# sheet from English to Malay. Only the header text changes; every data
# row keeps exactly the same values (the underlying shared-string table
# just gets re-packed as a natural side effect of removing the old
# header strings and appending the new ones).
#
#   A1 "Name"                -> "Nama"
#   D1 "Home Address"        -> "Alamat Rumah"
#   H1 "Place of Birth"      -> "Tempat Lahir"
#   J1 "Age"                 -> "Umur"
#   M1 "Parent Salary (RM)"  -> "Gaji"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Nama"
$ws.Range("D1").Value = "Alamat Rumah"
$ws.Range("J1").Value = "Umur"
$ws.Range("M1").Value = "Gaji"
$ws.Range("H1").Value = "Tempat Lahir"

# Best-fit-ish column widths to accommodate the new/changed header text
# (column C is untouched - it already has the right width for the date
# column and is not part of this edit).
$ws.Columns("A").ColumnWidth = 14.666666666666666
$ws.Columns("B").ColumnWidth = 13.666666666666666
$ws.Columns("D").ColumnWidth = 27.5

# Reset the view: scroll back to the top and move the selection to H1
# (previously the sheet was scrolled to row 13 with C21 selected).
$ws.Range("A1").Select()
$ws.Range("H1").Select()
